$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1924.2354
$ws.Range("I32").Value = 1910.8334
$ws.Range("K32").Value = 1910.8334
$ws.Range("M32").Value = -1584.8334
$ws.Range("H92").Value = 1250.65
$ws.Range("I92").Value = 1425.9375
$ws.Range("K92").Value = 1425.9375
$ws.Range("M92").Value = -177.9375
$ws.Range("H96").Value = 1497.96
$ws.Range("J96").Value = 2252.5
$ws.Range("L96").Value = 6757.5
$ws.Range("N96").Value = -9503.5
$ws.Range("H106").Value = 2081.139
$ws.Range("I106").Value = 3264.0557
$ws.Range("J106").Value = 898.2222
$ws.Range("K106").Value = 3264.0557
$ws.Range("L106").Value = 898.2222
$ws.Range("M106").Value = -2633.0557
$ws.Range("N106").Value = -2160.2222
$ws.Range("H132").Value = 7642
$ws.Range("I132").Value = 7642
$ws.Range("K132").Value = 22926
$ws.Range("M132").Value = -20396
$ws.Range("H135").Value = 4264.355
$ws.Range("I135").Value = 1121.0526
$ws.Range("J135").Value = 9241.25
$ws.Range("K135").Value = 10089.4734
$ws.Range("L135").Value = 83171.25
$ws.Range("M135").Value = -7554.473399999999
$ws.Range("N135").Value = -88241.25
$ws.Range("H138").Value = 3760.625
$ws.Range("I138").Value = 3489.6875
$ws.Range("J138").Value = 3869
$ws.Range("K138").Value = 10469.0625
$ws.Range("L138").Value = 11607
$ws.Range("M138").Value = -5329.0625
$ws.Range("N138").Value = -21887

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2859209.5
$ws.Range("I61").Value = 1678.1
$ws.Range("K61").Value = 1678.1
$ws.Range("M61").Value = -1466.1
$ws.Range("H88").Value = 2337.3076
$ws.Range("I88").Value = 1749.5
$ws.Range("J88").Value = 2598.5557
$ws.Range("K88").Value = 1749.5
$ws.Range("L88").Value = 2598.5557
$ws.Range("M88").Value = -1343.5
$ws.Range("N88").Value = -3410.5557
$ws.Range("H91").Value = 2337.3076
$ws.Range("I91").Value = 1749.5
$ws.Range("J91").Value = 2598.5557
$ws.Range("K91").Value = 1749.5
$ws.Range("L91").Value = 2598.5557
$ws.Range("M91").Value = -345.5
$ws.Range("N91").Value = -5406.5557
$ws.Range("H97").Value = 626.05
$ws.Range("I97").Value = 626.05
$ws.Range("K97").Value = 626.05
$ws.Range("M97").Value = -130.05
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H136").Value = 2859209.5
$ws.Range("I136").Value = 1678.1
$ws.Range("K136").Value = 5034.299999999999
$ws.Range("M136").Value = -2484.299999999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 589
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 664.6
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 664.6
$ws.Range("M22").Value = -227
$ws.Range("N22").Value = -1010.6
$ws.Range("H86").Value = 4632.1665
$ws.Range("I86").Value = 2257.6667
$ws.Range("J86").Value = 7006.6665
$ws.Range("K86").Value = 2257.6667
$ws.Range("L86").Value = 7006.6665
$ws.Range("M86").Value = -1134.6667
$ws.Range("N86").Value = -9252.666499999999
$ws.Range("H89").Value = 4632.1665
$ws.Range("I89").Value = 2257.6667
$ws.Range("J89").Value = 7006.6665
$ws.Range("K89").Value = 11288.3335
$ws.Range("L89").Value = 35033.3325
$ws.Range("M89").Value = -5672.333500000001
$ws.Range("N89").Value = -46265.3325

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 15000
$ws.Range("K39").Value = 15000
$ws.Range("M39").Value = -14609
$ws.Range("H49").Value = 15000
$ws.Range("I49").Value = 15000
$ws.Range("K49").Value = 15000
$ws.Range("M49").Value = -14818
$ws.Range("H98").Value = 80000
$ws.Range("J98").Value = 80000
$ws.Range("L98").Value = 80000
$ws.Range("N98").Value = -84492
$ws.Range("H99").Value = 2390.0667
$ws.Range("I99").Value = 2006.125
$ws.Range("J99").Value = 2828.8572
$ws.Range("K99").Value = 2006.125
$ws.Range("L99").Value = 2828.8572
$ws.Range("M99").Value = -508.125
$ws.Range("N99").Value = -5824.8572
$ws.Range("H120").Value = 51325
$ws.Range("J120").Value = 51325
$ws.Range("L120").Value = 51325
$ws.Range("N120").Value = -58583
$ws.Range("H126").Value = 2390.0667
$ws.Range("I126").Value = 2006.125
$ws.Range("J126").Value = 2828.8572
$ws.Range("K126").Value = 6018.375
$ws.Range("L126").Value = 8486.571599999999
$ws.Range("M126").Value = -3548.375
$ws.Range("N126").Value = -13426.5716
$ws.Range("H132").Value = 2559.4443
$ws.Range("I132").Value = 2559.4443
$ws.Range("K132").Value = 7678.3329
$ws.Range("M132").Value = -5148.3329
$ws.Range("H134").Value = 2450.7646
$ws.Range("I134").Value = 2551.6333
$ws.Range("J134").Value = 1694.25
$ws.Range("K134").Value = 7654.8999
$ws.Range("L134").Value = 5082.75
$ws.Range("M134").Value = -5119.8999
$ws.Range("N134").Value = -10152.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 7135.8276
$ws.Range("I64").Value = 4029.8
$ws.Range("J64").Value = 7782.9165
$ws.Range("K64").Value = 12089.4
$ws.Range("L64").Value = 23348.7495
$ws.Range("M64").Value = -11819.4
$ws.Range("N64").Value = -23888.7495
$ws.Range("H67").Value = 7135.8276
$ws.Range("I67").Value = 4029.8
$ws.Range("J67").Value = 7782.9165
$ws.Range("K67").Value = 12089.4
$ws.Range("L67").Value = 23348.7495
$ws.Range("M67").Value = -11153.4
$ws.Range("N67").Value = -25220.7495
$ws.Range("H107").Value = 1172
$ws.Range("J107").Value = 1316
$ws.Range("L107").Value = 3948
$ws.Range("N107").Value = -7788
$ws.Range("H119").Value = 8651
$ws.Range("I119").Value = 2971.8
$ws.Range("K119").Value = 8915.400000000001
$ws.Range("M119").Value = -4077.400000000001
$ws.Range("H141").Value = 5265
$ws.Range("I141").Value = 5265
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 15795
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -10615

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.86957
$ws.Range("I2").Value = 74.77778000000001
$ws.Range("J2").Value = 144.2
$ws.Range("K2").Value = 74.77778000000001
$ws.Range("L2").Value = 144.2
$ws.Range("M2").Value = 38.22221999999999
$ws.Range("N2").Value = -370.2
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 100
$ws.Range("K17").Value = 100
$ws.Range("M17").Value = 68
$ws.Range("H47").Value = 24000.5
$ws.Range("J47").Value = 24000.5
$ws.Range("L47").Value = 24000.5
$ws.Range("N47").Value = -25136.5
$ws.Range("H132").Value = 17562.666
$ws.Range("I132").Value = 10104.667
$ws.Range("K132").Value = 30314.001
$ws.Range("M132").Value = -27784.001

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 2000
$ws.Range("K4").Value = 2000
$ws.Range("M4").Value = -1887
$ws.Range("H22").Value = 3364.05
$ws.Range("I22").Value = 3143.6667
$ws.Range("J22").Value = 3544.3635
$ws.Range("K22").Value = 3143.6667
$ws.Range("L22").Value = 3544.3635
$ws.Range("M22").Value = -2848.6667
$ws.Range("N22").Value = -4134.363499999999
$ws.Range("H27").Value = 3364.05
$ws.Range("I27").Value = 3143.6667
$ws.Range("J27").Value = 3544.3635
$ws.Range("K27").Value = 3143.6667
$ws.Range("L27").Value = 3544.3635
$ws.Range("M27").Value = -3036.6667
$ws.Range("N27").Value = -3758.3635
$ws.Range("H28").Value = 2000
$ws.Range("I28").Value = 2000
$ws.Range("K28").Value = 2000
$ws.Range("M28").Value = -1768
$ws.Range("H30").Value = 1326.7778
$ws.Range("I30").Value = 871.8570999999999
$ws.Range("J30").Value = 2919
$ws.Range("K30").Value = 871.8570999999999
$ws.Range("L30").Value = 2919
$ws.Range("M30").Value = -763.8570999999999
$ws.Range("N30").Value = -3135
$ws.Range("H37").Value = 2000
$ws.Range("I37").Value = 2000
$ws.Range("K37").Value = 2000
$ws.Range("M37").Value = -1893
$ws.Range("H122").Value = 3891.389
$ws.Range("I122").Value = 3556.5
$ws.Range("J122").Value = 4310
$ws.Range("K122").Value = 10669.5
$ws.Range("L122").Value = 12930
$ws.Range("M122").Value = -8219.5
$ws.Range("N122").Value = -17830
$ws.Range("H132").Value = 3971489
$ws.Range("I132").Value = 7578666.5
$ws.Range("K132").Value = 22735999.5
$ws.Range("M132").Value = -22733469.5
$ws.Range("H139").Value = 86437.5
$ws.Range("J139").Value = 86437.5
$ws.Range("L139").Value = 86437.5
$ws.Range("N139").Value = -96717.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1722.5
$ws.Range("I113").Value = 1455.5714
$ws.Range("K113").Value = 4366.7142
$ws.Range("M113").Value = -2196.7142
$ws.Range("H132").Value = 55557556
$ws.Range("I132").Value = 55557556
$ws.Range("K132").Value = 166672668
$ws.Range("M132").Value = -166670138

# --- Removals ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N141").ClearContents()
